$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.475888013839722
$ws.Range("B1").Value = 2.992345094680786
$ws.Range("C1").Value = 5.216126918792725
$ws.Range("D1").Value = 0.6646848320960999
$ws.Range("E1").Value = 0.7860522270202637
